$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Summary")
$ws.Activate()

# Fill in the "Spint(42) - Day 8" block totals (Written / Execution / Review)
$ws.Range("C45").Value = 6936
$ws.Range("C46").Value = 2096
$ws.Range("C47").Value = 2096

# Move the active selection to H41 (cursor position when the file was saved)
$ws.Range("H41").Select()
